$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- B13 gets the "Good" cell style (matches style index 16 used elsewhere in the sheet) ---
$ws.Range("B13").Style = "Good"

# --- Clear the old summary row (row 38) before rebuilding the hours-spent block ---
$ws.Range("A38:C38").ClearContents()

# --- Row 26: header cell "hours spent" (previously lived in A27) ---
$ws.Range("A26").Value2 = "hours spent"

# --- Rows 27-37: sequence number in column A, hours value in column B ---
$ws.Range("A27").Value2 = 2
$ws.Range("B27").Value2 = 36

$ws.Range("A28").Value2 = 3
$ws.Range("B28").Value2 = 36

$ws.Range("A29").Value2 = 4
$ws.Range("B29").Value2 = 20

$ws.Range("A30").Value2 = 5
$ws.Range("B30").Value2 = 31

$ws.Range("A31").Value2 = 6
$ws.Range("B31").Value2 = 20

$ws.Range("A32").Value2 = 7
$ws.Range("B32").Value2 = 26

$ws.Range("A33").Value2 = 8
$ws.Range("B33").Value2 = 20

$ws.Range("A34").Value2 = 9
$ws.Range("B34").Value2 = 22

$ws.Range("A35").Value2 = 10
$ws.Range("B35").Value2 = 14

$ws.Range("A36").Value2 = 13
$ws.Range("B36").Value2 = 30

$ws.Range("A37").Value2 = 14
$ws.Range("B37").Formula = "=14+19"

# --- Row 42: new total row (replaces the old row 38 total) ---
$ws.Range("A42").Value2 = "Total hours spent so far"
$ws.Range("B42").Formula = "=SUM(B27:B40)"
$ws.Range("C42").Formula = "=B42*135"

# --- Sheet view: scroll so row 10 is at the top, matching topLeftCell="A10" ---
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Application.ActiveWindow.ScrollColumn = 1
